$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'72.270.61"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = '  +0.33%  '
$ws.Range("D3").Value = "'2.657.66"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = '  +1.69%  '
$ws.Range("E4").Value = '  -0.04%  '
$ws.Range("D5").Value = "'598.09"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  -0.62%  '
$ws.Range("D6").Value = "'175.12"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  -1.13%  '
$ws.Range("E7").Value = '  -0.02%  '
$ws.Range("D8").Value = "'0.524"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  +0.18%  '
$ws.Range("D9").Value = "'2.656.17"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  +1.66%  '
$ws.Range("E10").Value = '  -0.43%  '
$ws.Range("E11").Value = '  +2.51%  '
$ws.Range("D12").Value = "'0.356"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  +1.81%  '
$ws.Range("D13").Value = "'5.01"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  -0.14%  '
$ws.Range("D14").Value = "'3.149.21"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  +1.78%  '
$ws.Range("D15").Value = "'0.0000185"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  -0.53%  '
$ws.Range("D16").Value = "'72.106.43"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  +0.13%  '
$ws.Range("D17").Value = "'26.27"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  -0.35%  '
$ws.Range("D18").Value = "'2.659.67"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  +1.70%  '
$ws.Range("D19").Value = "'12.26"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  +6.66%  '
$ws.Range("D20").Value = "'8.23"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  +5.04%  '
$ws.Range("D21").Value = "'370.65"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  -2.57%  '
$ws.Range("D22").Value = "'4.17"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  +0.66%  '
$ws.Range("D23").Value = "'2.04"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  +3.16%  '
$ws.Range("D24").Value = "'72.00"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  -0.94%  '
$ws.Range("E25").Value = '  +0.24%  '
$ws.Range("D26").Value = "'4.33"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  -0.26%  '
$ws.Range("D27").Value = "'9.78"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  +0.01%  '
$ws.Range("D28").Value = "'2.793.48"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  +1.68%  '
$ws.Range("E29").Value = '  -0.20%  '
$ws.Range("D30").Value = "'0.0₃0969"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  +2.81%  '
$ws.Range("D31").Value = "'8.08"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  +1.32%  '
$ws.Range("D32").Value = "'501.40"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  -3.27%  '
$ws.Range("E33").Value = '  -1.43%  '
$ws.Range("E34").Value = '  +0.30%  '
$ws.Range("D35").Value = "'0.999"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  -0.06%  '
$ws.Range("D36").Value = "'163.15"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  -0.13%  '
$ws.Range("D37").Value = "'19.49"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  +1.65%  '
$ws.Range("B38").Value = 'Kaspa'
$ws.Range("C38").Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Range("D38").Value = "'0.112"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  +0.71%  '
$ws.Range("B39").Value = 'WhiteBITCoin'
$ws.Range("C39").Value = 'https://coinranking.com/coin/GE4c3_TbB+whitebitcoin-wbt'
$ws.Range("D39").Value = "'18.99"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  -0.53%  '
$ws.Range("E40").Value = '  -0.91%  '
$ws.Range("E41").Value = '  -2.38%  '
$ws.Range("E42").Value = '  +0.04%  '
$ws.Range("E43").Value = '  -0.22%  '
$ws.Range("E44").Value = '  +0.18%  '
$ws.Range("E45").Value = '  +0.59%  '
$ws.Range("B46").Value = 'Aave'
$ws.Range("C46").Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Range("D46").Value = "'156.18"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  +4.45%  '
$ws.Range("B47").Value = 'OKB'
$ws.Range("C47").Value = 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
$ws.Range("D47").Value = "'39.48"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  +0.09%  '
$ws.Range("D48").Value = "'0.559"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  +3.69%  '
$ws.Range("D49").Value = "'3.73"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  +1.78%  '
$ws.Range("D50").Value = "'1.73"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  +2.56%  '
$ws.Range("D51").Value = "'0.0755"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  -1.27%  '
